$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.46388812271717
$ws.Range("D2").Value = 10.19202271113368
$ws.Range("E2").Value = 13.99772134997894
$ws.Range("F2").Value = 28.00002925415852
$ws.Range("G2").Value = 27.9926249224201
$ws.Range("H2").Value = 13.21923212800053
$ws.Range("I2").Value = 42.00370031271751
$ws.Range("J2").Value = 9.644818864116688
$ws.Range("O2").Value = 20.1862790208076
$ws.Range("B3").Value = 14.62529716494162
$ws.Range("D3").Value = 10.1879277601034
$ws.Range("E3").Value = 13.96069026936255
$ws.Range("F3").Value = 27.72584099711391
$ws.Range("G3").Value = 27.34581590149833
$ws.Range("H3").Value = 13.19068893700772
$ws.Range("I3").Value = 39.58439637612544
$ws.Range("J3").Value = 9.646505256136761
$ws.Range("O3").Value = 20.02620259791259
$ws.Range("B4").Value = 14.08227189835594
$ws.Range("D4").Value = 10.18718855004983
$ws.Range("E4").Value = 13.94095794876877
$ws.Range("F4").Value = 27.56376897813316
$ws.Range("G4").Value = 26.9496947766748
$ws.Range("H4").Value = 13.17586455310832
$ws.Range("I4").Value = 38.02187358361297
$ws.Range("J4").Value = 9.649197674905507
$ws.Range("O4").Value = 19.93251056327971
$ws.Range("B5").Value = 13.85402509196118
$ws.Range("D5").Value = 10.18733472167911
$ws.Range("E5").Value = 13.93367796296435
$ws.Range("F5").Value = 27.49937375959707
$ws.Range("G5").Value = 26.78879404821963
$ws.Range("H5").Value = 13.17050689977642
$ws.Range("I5").Value = 37.36626242888532
$ws.Range("J5").Value = 9.650711919510597
$ws.Range("O5").Value = 19.89552268896703
$ws.Range("B6").Value = 13.81570835831321
$ws.Range("D6").Value = 10.18738604175381
$ws.Range("E6").Value = 13.93251524037363
$ws.Range("F6").Value = 27.488782674454
$ws.Range("G6").Value = 26.76211597885145
$ws.Range("H6").Value = 13.16965865087502
$ws.Range("I6").Value = 37.25627603765865
$ws.Range("J6").Value = 9.65098855693036
$ws.Range("O6").Value = 19.8894539280562
$ws.Range("B7").Value = 14.07922169308636
$ws.Range("D7").Value = 10.18718870862609
$ws.Range("E7").Value = 13.94085668000986
$ws.Range("F7").Value = 27.56289374781946
$ws.Range("G7").Value = 26.94752234594259
$ws.Range("H7").Value = 13.17578952584707
$ws.Range("I7").Value = 38.01310743513349
$ws.Range("J7").Value = 9.649216407453391
$ws.Range("O7").Value = 19.93200685707592
$ws.Range("B8").Value = 15.18063999774328
$ws.Range("D8").Value = 10.19024303432427
$ws.Range("E8").Value = 13.98433288131076
$ws.Range("F8").Value = 27.90422929772455
$ws.Range("G8").Value = 27.76954980634348
$ws.Range("H8").Value = 13.20883245916042
$ws.Range("I8").Value = 41.18578900812587
$ws.Range("J8").Value = 9.645056714851505
$ws.Range("O8").Value = 20.13015363980903
$ws.Range("B9").Value = 17.11354650569201
$ws.Range("D9").Value = 10.2102635689738
$ws.Range("E9").Value = 14.09315794686801
$ws.Range("F9").Value = 28.6201392906549
$ws.Range("G9").Value = 29.37840297970911
$ws.Range("H9").Value = 13.29485483016452
$ws.Range("I9").Value = 46.77879708202151
$ws.Range("J9").Value = 9.650023774114835
$ws.Range("O9").Value = 20.55345977065718
$ws.Range("B10").Value = 18.39153772040988
$ws.Range("D10").Value = 10.23344268846874
$ws.Range("E10").Value = 14.1870649376271
$ws.Range("F10").Value = 29.16991955979584
$ws.Range("G10").Value = 30.54415329235382
$ws.Range("H10").Value = 13.37067982789178
$ws.Range("I10").Value = 50.48748964831373
$ws.Range("J10").Value = 9.661635036204377
$ws.Range("O10").Value = 20.88320280852609
$ws.Range("B11").Value = 18.94157951104755
$ws.Range("D11").Value = 10.24580427003773
$ws.Range("E11").Value = 14.23271206016586
$ws.Range("F11").Value = 29.42418676710619
$ws.Range("G11").Value = 31.06815545768676
$ws.Range("H11").Value = 13.40783664559342
$ws.Range("I11").Value = 52.08534703260537
$ws.Range("J11").Value = 9.66863450946393
$ws.Range("O11").Value = 21.036730358788
$ws.Range("B12").Value = 19.14533481163442
$ws.Range("D12").Value = 10.25074430124067
$ws.Range("E12").Value = 14.25040869153869
$ws.Range("F12").Value = 29.52098477406991
$ws.Range("G12").Value = 31.26546135151553
$ws.Range("H12").Value = 13.42228238153746
$ws.Range("I12").Value = 52.67744494454099
$ws.Range("J12").Value = 9.671530688125269
$ws.Range("O12").Value = 21.09532649437566
$ws.Range("B13").Value = 19.10165446781605
$ws.Range("D13").Value = 10.2496688993483
$ws.Range("E13").Value = 14.24657929015571
$ws.Range("F13").Value = 29.50011618201164
$ws.Range("G13").Value = 31.22302119215394
$ws.Range("H13").Value = 13.4191546736729
$ws.Range("I13").Value = 52.55050504930303
$ws.Range("J13").Value = 9.670896042929174
$ws.Range("O13").Value = 21.08268715005322
$ws.Range("B14").Value = 18.95843357864365
$ws.Range("D14").Value = 10.24620551355724
$ws.Range("E14").Value = 14.23415979427644
$ws.Range("F14").Value = 29.43214062403234
$ws.Range("G14").Value = 31.08441160933745
$ws.Range("H14").Value = 13.40901763072708
$ws.Range("I14").Value = 52.13432002565632
$ws.Range("J14").Value = 9.668867865316763
$ws.Range("O14").Value = 21.041542187005
$ws.Range("B15").Value = 18.87011552320599
$ws.Range("D15").Value = 10.24411774134096
$ws.Range("E15").Value = 14.22660571722196
$ws.Range("F15").Value = 29.3905677953432
$ws.Range("G15").Value = 30.99935692190976
$ws.Range("H15").Value = 13.40285702873039
$ws.Range("I15").Value = 51.87770123785431
$ws.Range("J15").Value = 9.667657493565205
$ws.Range("O15").Value = 21.01639798227277
$ws.Range("B16").Value = 18.3549581758265
$ws.Range("D16").Value = 10.23267120236419
$ws.Range("E16").Value = 14.18413987375155
$ws.Range("F16").Value = 29.15337906325784
$ws.Range("G16").Value = 30.50976287739119
$ws.Range("H16").Value = 13.36830448962632
$ws.Range("I16").Value = 50.38125698588405
$ws.Range("J16").Value = 9.66121205336599
$ws.Range("O16").Value = 20.87323624657447
$ws.Range("B17").Value = 18.03087938850338
$ws.Range("D17").Value = 10.22611296863152
$ws.Range("E17").Value = 14.15883154440653
$ws.Range("F17").Value = 29.00887749292154
$ws.Range("G17").Value = 30.20763995618443
$ws.Range("H17").Value = 13.34778449073917
$ws.Range("I17").Value = 49.44026495104949
$ws.Range("J17").Value = 9.657696969306752
$ws.Range("O17").Value = 20.78627997669675
$ws.Range("B18").Value = 17.84153069387959
$ws.Range("D18").Value = 10.22251205291303
$ws.Range("E18").Value = 14.14455113750949
$ws.Range("F18").Value = 28.92616169188516
$ws.Range("G18").Value = 30.03328550724626
$ws.Range("H18").Value = 13.33623314540066
$ws.Range("I18").Value = 48.89063678530206
$ws.Range("J18").Value = 9.655836931840525
$ws.Range("O18").Value = 20.73659984591648
$ws.Range("B19").Value = 17.77691535278072
$ws.Range("D19").Value = 10.22132231765938
$ws.Range("E19").Value = 14.1397637826668
$ws.Range("F19").Value = 28.89822649347399
$ws.Range("G19").Value = 29.97415906526003
$ws.Range("H19").Value = 13.33236543278871
$ws.Range("I19").Value = 48.7031053676155
$ws.Range("J19").Value = 9.655234972994768
$ws.Range("O19").Value = 20.719837975887
$ws.Range("B20").Value = 18.06568338928438
$ws.Range("D20").Value = 10.22679340089665
$ws.Range("E20").Value = 14.16149713654993
$ws.Range("F20").Value = 29.02421940351402
$ws.Range("G20").Value = 30.23986323467777
$ws.Range("H20").Value = 13.34994293450564
$ws.Range("I20").Value = 49.54130510666898
$ws.Range("J20").Value = 9.658054425068856
$ws.Range("O20").Value = 20.79550231663106
$ws.Range("B21").Value = 19.00062425561376
$ws.Range("D21").Value = 10.24721578547861
$ws.Range("E21").Value = 14.23779662850087
$ws.Range("F21").Value = 29.45209347683171
$ws.Range("G21").Value = 31.12515671041986
$ws.Range("H21").Value = 13.41198500623147
$ws.Range("I21").Value = 52.25691677039583
$ws.Range("J21").Value = 9.669456935443923
$ws.Range("O21").Value = 21.05361540318239
$ws.Range("B22").Value = 19.58522931221873
$ws.Range("D22").Value = 10.26207120930413
$ws.Range("E22").Value = 14.29005360826744
$ws.Range("F22").Value = 29.73468108382609
$ws.Range("G22").Value = 31.697115733862
$ws.Range("H22").Value = 13.4547163867481
$ws.Range("I22").Value = 53.9560509119212
$ws.Range("J22").Value = 9.678340111134078
$ws.Range("O22").Value = 21.22495736481463
$ws.Range("B23").Value = 19.27564066207565
$ws.Range("D23").Value = 10.25400542469682
$ws.Range("E23").Value = 14.26194779234118
$ws.Range("F23").Value = 29.58361782606772
$ws.Range("G23").Value = 31.39252465651737
$ws.Range("H23").Value = 13.43171277708358
$ws.Range("I23").Value = 53.05615300520156
$ws.Range("J23").Value = 9.673468545733964
$ws.Range("O23").Value = 21.13328259245116
$ws.Range("B24").Value = 18.04995794057753
$ws.Range("D24").Value = 10.22648524933143
$ws.Range("E24").Value = 14.16029118155925
$ws.Range("F24").Value = 29.01728219750559
$ws.Range("G24").Value = 30.22529712948352
$ws.Range("H24").Value = 13.34896633560851
$ws.Range("I24").Value = 49.49565175058866
$ws.Range("J24").Value = 9.657892318257417
$ws.Range("O24").Value = 20.79133192062994
$ws.Range("B25").Value = 16.61539234780995
$ws.Range("D25").Value = 10.20335309967254
$ws.Range("E25").Value = 14.0612345977308
$ws.Range("F25").Value = 28.42194760795592
$ws.Range("G25").Value = 28.94502986164266
$ws.Range("H25").Value = 13.26934136391781
$ws.Range("I25").Value = 45.33509062920486
$ws.Range("J25").Value = 9.647278360295349
$ws.Range("O25").Value = 20.43547031540893
